# Applies the "first values to Singapore archetypes database" commit:
# adds the first data values (T0/T1/T2/T6/T7 construction-type codes and
# a couple of numeric ratios) to the ARCHITECTURE and HVAC sheets of the
# Singapore (SIN) construction_properties archetypes workbook.
#
# The order in which the values below are written matters: it reproduces
# the order in which the new strings were first introduced into the
# workbook's shared-string table by the original author.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ARCHITECTURE")
$ws2 = $wb.Worksheets.Item("HVAC")

# ---------------------------------------------------------------------
# ARCHITECTURE sheet
# ---------------------------------------------------------------------

# Row 5 (OFFICE): type_leak
$ws1.Range("H5").Value = "T1"

# type_cons (col G): every data row (2-19) gets "T2"
$ws1.Range("G2:G19").Value = "T2"

# type_roof (col J): every data row (2-19) gets "T7"
$ws1.Range("J2:J19").Value = "T7"

# type_wall (col K): most rows get "T6", the first two (MULTI_RES,
# SINGLE_RES) get "T2"; some rows (HOTEL, FOODSTORE, INDUSTRIAL,
# HOSPITAL, SWIMMING, SERVERROOM, PARKING, COOLROOM) are left blank
$ws1.Range("K2:K3").Value = "T2"
$ws1.Range("K5:K6").Value = "T6"
$ws1.Range("K8").Value = "T6"
$ws1.Range("K10").Value = "T6"
$ws1.Range("K12").Value = "T6"
$ws1.Range("K17:K19").Value = "T6"

# Row 5 (OFFICE): type_win
$ws1.Range("I5").Value = "T2"

# Row 5 (OFFICE) and row 17 (LAB): win_wall ratio
$ws1.Range("F5").Value = 0.35
$ws1.Range("F17").Value = 0.11

# ---------------------------------------------------------------------
# HVAC sheet
# ---------------------------------------------------------------------

# type_hs (col E): every data row (2-19) gets "T0"
$ws2.Range("E2:E19").Value = "T0"

# ---------------------------------------------------------------------
# Restore on-screen selections / active sheet to match the saved state
# ---------------------------------------------------------------------

$ws2.Activate() | Out-Null
$ws2.Range("E2:E19").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("G19").Select() | Out-Null

Write-Host "Applied Singapore archetype values"
